{"js": "// Update the worksheet date and all multiplication answers to the new set\n// of values, matching the OOXML diff exactly. Every old text value is\n// unique in the document, so a simple search-and-replace per pair is safe.\nconst replacements = [\n  [\"2025-11-03 Monday\", \"2025-11-04 Tuesday\"],\n  [\"57\u00d714=798\", \"44\u00d763=2772\"],\n  [\"93\u00d757=5301\", \"83\u00d727=2241\"],\n  [\"62\u00d749=3038\", \"47\u00d764=3008\"],\n  [\"89\u00d749=4361\", \"57\u00d751=2907\"],\n  [\"33\u00d777=2541\", \"44\u00d739=1716\"],\n  [\"34\u00d798=3332\", \"34\u00d732=1088\"],\n  [\"39\u00d750=1950\", \"91\u00d788=8008\"],\n  [\"18\u00d788=1584\", \"41\u00d791=3731\"],\n  [\"20\u00d740=800\", \"75\u00d771=5325\"],\n  [\"55\u00d758=3190\", \"36\u00d776=2736\"],\n  [\"67\u00d775=5025\", \"32\u00d758=1856\"],\n  [\"93\u00d740=3720\", \"76\u00d782=6232\"],\n  [\"73\u00d754=3942\", \"19\u00d721=399\"],\n  [\"77\u00d750=3850\", \"59\u00d782=4838\"],\n  [\"97\u00d743=4171\", \"60\u00d730=1800\"],\n  [\"52\u00d790=4680\", \"46\u00d797=4462\"],\n  [\"49\u00d731=1519\", \"97\u00d751=4947\"],\n  [\"44\u00d735=1540\", \"37\u00d728=1036\"],\n  [\"38\u00d764=2432\", \"36\u00d725=900\"],\n  [\"64\u00d743=2752\", \"11\u00d768=748\"],\n  [\"30\u00d789=2670\", \"41\u00d732=1312\"],\n  [\"95\u00d761=5795\", \"70\u00d737=2590\"],\n  [\"16\u00d733=528\", \"46\u00d764=2944\"],\n  [\"67\u00d727=1809\", \"49\u00d774=3626\"],\n  [\"21\u00d720=420\", \"36\u00d795=3420\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and all multiplication answers to the new set\n# of values, matching the OOXML diff exactly. Every old text value is\n# unique in the document, so Find/Replace (wdReplaceAll) per pair is safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-11-03 Monday\", \"2025-11-04 Tuesday\"),\n    @(\"57\u00d714=798\", \"44\u00d763=2772\"),\n    @(\"93\u00d757=5301\", \"83\u00d727=2241\"),\n    @(\"62\u00d749=3038\", \"47\u00d764=3008\"),\n    @(\"89\u00d749=4361\", \"57\u00d751=2907\"),\n    @(\"33\u00d777=2541\", \"44\u00d739=1716\"),\n    @(\"34\u00d798=3332\", \"34\u00d732=1088\"),\n    @(\"39\u00d750=1950\", \"91\u00d788=8008\"),\n    @(\"18\u00d788=1584\", \"41\u00d791=3731\"),\n    @(\"20\u00d740=800\", \"75\u00d771=5325\"),\n    @(\"55\u00d758=3190\", \"36\u00d776=2736\"),\n    @(\"67\u00d775=5025\", \"32\u00d758=1856\"),\n    @(\"93\u00d740=3720\", \"76\u00d782=6232\"),\n    @(\"73\u00d754=3942\", \"19\u00d721=399\"),\n    @(\"77\u00d750=3850\", \"59\u00d782=4838\"),\n    @(\"97\u00d743=4171\", \"60\u00d730=1800\"),\n    @(\"52\u00d790=4680\", \"46\u00d797=4462\"),\n    @(\"49\u00d731=1519\", \"97\u00d751=4947\"),\n    @(\"44\u00d735=1540\", \"37\u00d728=1036\"),\n    @(\"38\u00d764=2432\", \"36\u00d725=900\"),\n    @(\"64\u00d743=2752\", \"11\u00d768=748\"),\n    @(\"30\u00d789=2670\", \"41\u00d732=1312\"),\n    @(\"95\u00d761=5795\", \"70\u00d737=2590\"),\n    @(\"16\u00d733=528\", \"46\u00d764=2944\"),\n    @(\"67\u00d727=1809\", \"49\u00d774=3626\"),\n    @(\"21\u00d720=420\", \"36\u00d795=3420\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\n$d.Save()\n"}
